$wb = $excel.ActiveWorkbook

# ----- Sheet ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 121
$ws.Range("H121").Value = 1072.7084
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 1072.7084
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 3218.1252
$ws.Range("M121").ClearContents()
$ws.Range("N121").Value = -6712.1252
# Row 129
$ws.Range("H129").Value = 752.125
$ws.Range("I129").Value = 470.44446
$ws.Range("K129").Value = 1411.33338
$ws.Range("M129").Value = 3588.66662
# Row 138
$ws.Range("H138").Value = 2314.2886
$ws.Range("I138").Value = 2280.7144
$ws.Range("J138").Value = 2326.658
$ws.Range("K138").Value = 6842.1432
$ws.Range("L138").Value = 6979.974
$ws.Range("M138").Value = -1702.1432
$ws.Range("N138").Value = -17259.974

# ----- Sheet ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 24
$ws.Range("H24").Value = 29199.545
$ws.Range("J24").Value = 29199.545
$ws.Range("L24").Value = 29199.545
$ws.Range("N24").Value = -29947.545
# Row 28
$ws.Range("H28").Value = 6480
$ws.Range("I28").Value = 6480
$ws.Range("K28").Value = 6480
$ws.Range("M28").Value = -6288
# Row 82
$ws.Range("H82").Value = 20600
$ws.Range("J82").Value = 20250
$ws.Range("L82").Value = 20250
$ws.Range("N82").Value = -20972
# Row 85
$ws.Range("H85").Value = 20600
$ws.Range("J85").Value = 20250
$ws.Range("L85").Value = 20250
$ws.Range("N85").Value = -22746
# Row 99
$ws.Range("H99").Value = 6480
$ws.Range("I99").Value = 6480
$ws.Range("K99").Value = 6480
$ws.Range("M99").Value = -3485
# Row 100
$ws.Range("H100").Value = 29199.545
$ws.Range("J100").Value = 29199.545
$ws.Range("L100").Value = 29199.545
$ws.Range("N100").Value = -31363.545
# Row 122
$ws.Range("H122").Value = 1613.875
$ws.Range("I122").Value = 2344.4
$ws.Range("J122").Value = 1281.8182
$ws.Range("K122").Value = 7033.200000000001
$ws.Range("L122").Value = 3845.4546
$ws.Range("M122").Value = -4583.200000000001
$ws.Range("N122").Value = -8745.454600000001
# Row 125
$ws.Range("H125").Value = 150044940
$ws.Range("J125").Value = 150044940
$ws.Range("L125").Value = 150044940
$ws.Range("N125").Value = -150054780
# Row 129
$ws.Range("H129").Value = 48699.5
$ws.Range("J129").Value = 48699.5
$ws.Range("L129").Value = 48699.5
$ws.Range("N129").Value = -58699.5

# ----- Sheet CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 87
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
# Row 90
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
# Row 100
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

# ----- Sheet CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 76
$ws.Range("H76").Value = 5707.3076
$ws.Range("I76").Value = 1998.3334
$ws.Range("J76").Value = 6820
$ws.Range("K76").Value = 5995.0002
$ws.Range("L76").Value = 20460
$ws.Range("M76").Value = -5612.0002
$ws.Range("N76").Value = -21226
# Row 79
$ws.Range("H79").Value = 5707.3076
$ws.Range("I79").Value = 1998.3334
$ws.Range("J79").Value = 6820
$ws.Range("K79").Value = 5995.0002
$ws.Range("L79").Value = 20460
$ws.Range("M79").Value = -4669.0002
$ws.Range("N79").Value = -23112
# Row 131
$ws.Range("H131").Value = 880.202
$ws.Range("I131").Value = 550
$ws.Range("J131").Value = 921.4773
$ws.Range("K131").Value = 1650
$ws.Range("L131").Value = 2764.4319
$ws.Range("M131").Value = 3390
$ws.Range("N131").Value = -12844.4319

# ----- Sheet GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 1382.2084
$ws.Range("I113").Value = 1295.2222
$ws.Range("J113").Value = 1643.1666
$ws.Range("K113").Value = 1295.2222
$ws.Range("L113").Value = 1643.1666
$ws.Range("M113").Value = 874.7778000000001
$ws.Range("N113").Value = -5983.1666
# Row 122
$ws.Range("H122").Value = 2040.5
$ws.Range("I122").Value = 2067.5
$ws.Range("K122").Value = 6202.5
$ws.Range("M122").Value = -3752.5
# Row 124
$ws.Range("H124").Value = 56500
$ws.Range("J124").Value = 56500
$ws.Range("L124").Value = 56500
$ws.Range("N124").Value = -66320

# ----- Sheet LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2821.7778
$ws.Range("I7").Value = 2416
$ws.Range("J7").Value = 3633.3333
$ws.Range("K7").Value = 2416
$ws.Range("L7").Value = 3633.3333
$ws.Range("M7").Value = -2304
$ws.Range("N7").Value = -3857.3333
# Row 9
$ws.Range("H9").Value = 5454.143
$ws.Range("I9").Value = 529.8333
$ws.Range("J9").Value = 35000
$ws.Range("K9").Value = 529.8333
$ws.Range("L9").Value = 35000
$ws.Range("M9").Value = -305.8333
$ws.Range("N9").Value = -35448
# Row 109
$ws.Range("H109").Value = 20921.666
$ws.Range("J109").Value = 20921.666
$ws.Range("L109").Value = 20921.666
$ws.Range("N109").Value = -23695.666
# Row 118
$ws.Range("H118").Value = 34562.668
$ws.Range("J118").Value = 34562.668
$ws.Range("L118").Value = 34562.668
$ws.Range("N118").Value = -37876.668
# Row 122
$ws.Range("H122").Value = 3934
$ws.Range("I122").Value = 2500
$ws.Range("J122").Value = 4548.5713
$ws.Range("K122").Value = 7500
$ws.Range("L122").Value = 13645.7139
$ws.Range("M122").Value = -5050
$ws.Range("N122").Value = -18545.7139
# Row 126
$ws.Range("H126").Value = 2821.7778
$ws.Range("I126").Value = 2416
$ws.Range("J126").Value = 3633.3333
$ws.Range("K126").Value = 7248
$ws.Range("L126").Value = 10899.9999
$ws.Range("M126").Value = -4778
$ws.Range("N126").Value = -15839.9999
# Row 127
$ws.Range("H127").Value = 55290
$ws.Range("J127").Value = 55290
$ws.Range("L127").Value = 55290
$ws.Range("N127").Value = -65210

# ----- Sheet WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 75
$ws.Range("H75").Value = 29875
$ws.Range("J75").Value = 29875
$ws.Range("L75").Value = 29875
$ws.Range("N75").Value = -31747
# Row 78
$ws.Range("H78").Value = 29875
$ws.Range("J78").Value = 29875
$ws.Range("L78").Value = 89625
$ws.Range("N78").Value = -98985
# Row 107
$ws.Range("H107").Value = 1768.7273
$ws.Range("I107").Value = 1084.0952
$ws.Range("J107").Value = 2966.8333
$ws.Range("K107").Value = 3252.2856
$ws.Range("L107").Value = 8900.499899999999
$ws.Range("M107").Value = -1332.2856
$ws.Range("N107").Value = -12740.4999
# Row 121
$ws.Range("H121").Value = 29600
$ws.Range("J121").Value = 29600
$ws.Range("L121").Value = 29600
$ws.Range("N121").Value = -33094
# Row 122
$ws.Range("H122").Value = 2132.2222
$ws.Range("I122").Value = 1875.7273
$ws.Range("K122").Value = 5627.1819
$ws.Range("M122").Value = -3177.1819

Write-Output "Applied all profit-sheet updates"